$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 21: new "Họp dự án" entry (17/9/2014) ----
# Copy the formatting of the existing, style-identical "Họp dự án" row (15)
# onto the new row 21 before writing values.
$ws.Range("A15:K15").Copy()
$ws.Range("A21:K21").PasteSpecial(-4122)

$ws.Range("A21").Value = "17/9/2014"
$ws.Range("B21").Value = "17/9/2014"
$ws.Range("C21").Value = "17/9/2014"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = "17/9/2014"
$ws.Range("F21").Value = "Họp dự án"
$ws.Range("G21").Value = "*Tổng hợp kiến thức cá nhân mỗi người`n*Kế hoạch tiếp theo"
$ws.Range("H21").Value = "OK"
$ws.Range("I21").Value = "OK"
$ws.Range("J21").Value = "OK"
$ws.Range("K21").Value = "OK"

$ws.Rows("21").RowHeight = 30

# ---- Row 24: new "Xác nhận mail" task entry (27/9/2014) ----
# Copy per-cell formatting from row 8 (same column pattern: A,B,C,D,F,G,I)
# so we don't introduce E/H/J/K cells that row 8 (as a whole) would add.
$ws.Range("A8").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F24").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("G24").PasteSpecial(-4122)
$ws.Range("I8").Copy()
$ws.Range("I24").PasteSpecial(-4122)

$ws.Range("A24").Value = "17/9/2014"
$ws.Range("F24").Value = "Xác nhận mail"
$ws.Range("G24").Value = "*Tìm hiểu xác nhận mail đăng nhập`n"
$ws.Range("I24").Value = "`n*Todo: Tìm hiểu về cơ chế sinh link, key để đăng ký acc cho website"
$ws.Range("B24").Value = "27/9/2014"

$ws.Rows("24").RowHeight = 30

# ---- View state: selection moved to B24, scrolled down a bit ----
$ws.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
